# Trade #3 closed at 2026-02-16 22:52:14 - base_strategy DOWN +0.000%
# Append a new row (row 4) with the third trade to both the "All Trades"
# and "base_strategy" sheets (they mirror each other).

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column B ("Date") holds a literal "yyyy-mm-dd" string in this sheet,
    # not a real date. Writing that text straight into .Value lets Excel's
    # autodetection reinterpret it as a date serial number, so instead
    # copy the identical date string already sitting in B3 down into B4 -
    # a copy/paste keeps the original literal-text representation intact.
    $ws.Range("B3").Copy()
    $ws.Range("B4").PasteSpecial()

    $ws.Cells.Item(4, 1).Value = 3
    $ws.Cells.Item(4, 3).Value = "22:52:14"
    $ws.Cells.Item(4, 4).Value = "base_strategy"
    $ws.Cells.Item(4, 5).Value = "DOWN"
    $ws.Cells.Item(4, 6).Value = 49.999998
    $ws.Cells.Item(4, 8).Value = "OPEN"
    $ws.Cells.Item(4, 9).Value = 0
    $ws.Cells.Item(4, 10).Value = 0
    $ws.Cells.Item(4, 11).Value = 100
    $ws.Cells.Item(4, 12).Value = 0
    $ws.Cells.Item(4, 13).Value = 0
    $ws.Cells.Item(4, 14).Value = 0.6
    $ws.Cells.Item(4, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(4, 17).Value = 0
}
